# Auto-committed on 2023/04/06 週四  9:52:13.34
# Adds two new trailing header columns ("案件隸屬單位" / "企金別") to the
# "放款餘額明細表" worksheet, widens the new column, and updates the
# active selection of the frozen (bottom-right) pane.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the two new header cells in row 1 (columns Y and Z) -----------
$ws.Cells.Item(1, 25).Value = "案件隸屬單位"
$ws.Cells.Item(1, 26).Value = "企金別"

# --- Give the new "案件隸屬單位" column (Y / column 25) a sensible width --
$ws.Columns.Item(25).ColumnWidth = 14

# --- Update the active selection within the frozen bottom-right pane ---
$ws.Range("X14").Select() | Out-Null
